$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-06-20 Friday" "2025-06-21 Saturday"

Replace-Text "122×4=488" "323×8=2584"
Replace-Text "385×5=1925" "372×4=1488"
Replace-Text "766×6=4596" "730×7=5110"
Replace-Text "140×7=980" "977×9=8793"
Replace-Text "810×6=4860" "910×4=3640"

Replace-Text "757×8=6056" "760×4=3040"
Replace-Text "608×8=4864" "680×9=6120"
Replace-Text "308×2=616" "960×6=5760"
Replace-Text "991×3=2973" "579×2=1158"
Replace-Text "883×9=7947" "256×8=2048"

Replace-Text "951×2=1902" "630×9=5670"
Replace-Text "244×2=488" "966×6=5796"
Replace-Text "501×7=3507" "299×6=1794"
Replace-Text "828×9=7452" "828×5=4140"
Replace-Text "969×2=1938" "170×5=850"

Replace-Text "329×4=1316" "735×6=4410"
Replace-Text "141×3=423" "579×6=3474"
Replace-Text "751×8=6008" "259×6=1554"
Replace-Text "771×2=1542" "727×8=5816"
Replace-Text "808×5=4040" "337×7=2359"

Replace-Text "754×3=2262" "526×4=2104"
Replace-Text "179×6=1074" "901×7=6307"
Replace-Text "443×8=3544" "966×6=5796"
Replace-Text "916×7=6412" "671×3=2013"
Replace-Text "657×8=5256" "780×6=4680"
